$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of B2 (previously held the value 23253374414)
$ws.Range("B2").ClearContents()

# Update the active selection to B8 to match the saved view state
$ws.Range("B8").Select()
